$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: the footnote block currently sitting on row 14 moves down
# two rows (to row 16) to leave space for the new "Validação/Total" summary
# row that lands on row 13.
$ws.Rows("14:15").Insert()
$ws.Rows("16").EntireRow.AutoFit()

# --- Header row (row 3): two new trailing columns (Validação, Total) and
# the Brier column gets an explicit footnote marker.
$ws.Range("G3").Value = "Validação"
$ws.Range("H3").Value = "Total"
$ws.Range("C3").Value = "Brier *"

# --- Existing data rows (5-12): add "Validação" (G) and "Total" (H) columns
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 9101

$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 9101

$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 9101

$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 9101

$ws.Range("E9").Value = 5824
$ws.Range("F9").Value = 1821
$ws.Range("G9").Value = 1456
$ws.Range("H9").Value = 9101

$ws.Range("E10").Value = 5824
$ws.Range("F10").Value = 1821
$ws.Range("G10").Value = 1456
$ws.Range("H10").Value = 9101

$ws.Range("E11").Value = 5824
$ws.Range("F11").Value = 1821
$ws.Range("G11").Value = 1456
$ws.Range("H11").Value = 9101

$ws.Range("E12").Value = 5824
$ws.Range("F12").Value = 1821
$ws.Range("G12").Value = 1456
$ws.Range("H12").Value = 9101

# --- Row 13: percentage split of the Treino/Teste/Validação/Total columns
$ws.Range("E13").Value = 0.64
$ws.Range("F13").Value = 0.2
$ws.Range("G13").Value = 0.16
$ws.Range("H13").Value = 1
$ws.Range("E13:H13").NumberFormat = "0%"

# --- Restore the selection to match the relocated footnote block
$ws.Range("C16:E16").Select()
